$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = "39,87 TRY - 79,76 TRY - 797,68 TRY"
$ws.Range("G4").Value = "27,84 TRY - 55,69 TRY - 398,83 TRY"
$ws.Range("G5").Value = "7,97 TRY - 15,96 TRY - 199,41 TRY"

$ws.Range("C6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("G6").Value = "8.300,01 TL - 99,71 TL"
$ws.Range("J6").Value = "8.300,01 TL - 199,41 TL"

$ws.Range("G8").Value = "19,94 TRY - 39,88 TRY - 398,84 TRY"
$ws.Range("G9").Value = "13,92 TRY - 27,85 TRY - 199,42 TRY"
$ws.Range("G10").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"
$ws.Range("G11").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"

$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("G12").Value = "Şube (Kasadan): %0,5; Şube (Hesaptan): %0,75; İnternet: 15 USD"

$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("F13").Value = ""
$ws.Range("J13").Value = "Hesaba: Asgari 1 TL | Azami 995,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

$ws.Range("C14").Value = "40.000 TL - 2.485,72 TL"
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = "8.300 TL - 7,97 TL"
$ws.Range("J14").Value = "1.554,97 TL - 7.784 TL"

$ws.Range("J15").Value = "%0,5 Asgari Tutar: 375 TL Azami Tutar: 375 TL / 6.500 TL"
$ws.Range("J17").Value = "%0,5 Asgari Tutar: 350 TL Azami Tutar: 350 TL / 5.500 TL"
$ws.Range("J20").Value = "200 TL"
$ws.Range("J21").Value = "%0,5 Asgari Tutar: 500 TL Azami Tutar: 500 TL / 5.000 TL"
$ws.Range("J22").Value = "%0,5 Asgari Tutar: 1.500 TL Azami Tutar: 1.500 TL / 9.000 TL"
$ws.Range("J23").Value = "75 TL"
$ws.Range("J24").Value = "500 TL"
$ws.Range("J25").Value = "750 TL"
